$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header / summary block updates
# ---------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"
$ws.Range("C8").Value = 4633.38
$ws.Range("C9").Value = 35
$ws.Range("G10").Value = ""

# ---------------------------------------------------------------------
# 2) Remove the obsolete "PLA-TAG" line item (old row 17). Deleting the
#    entire row shifts every subsequent row up by one and auto-adjusts
#    the merged TOTAL range.
# ---------------------------------------------------------------------
$ws.Rows(17).Delete()

# The row delete carries each row's original formatting along with its
# content, but the sheet's zebra-striping is meant to stay keyed to row
# *position* (odd/even), not to the content that happens to occupy it.
# Re-stamp the alternating two-row format pattern across the whole
# line-item block so the banding lines back up.
$ws.Range("A16:H17").Copy()
$ws.Range("A16:H51").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Populate the previously-zeroed "Pricing" column (H) for every
#    remaining line-item row (now rows 16-50) plus the TOTAL row (51).
# ---------------------------------------------------------------------
$pricing = @{
    16 = 31.72
    17 = 286.25
    18 = 95.15000000000001
    19 = 40.58
    20 = 7.38
    21 = 128.82
    22 = 118.38
    23 = 648.53
    24 = 61.51
    25 = 286.25
    26 = 95.15000000000001
    27 = 195.83
    28 = 60.87
    29 = 128.82
    30 = 128.82
    31 = 118.38
    32 = 118.38
    33 = 81.16
    34 = 14.76
    35 = 195.83
    36 = 286.25
    37 = 95.15000000000001
    38 = 128.82
    39 = 118.38
    40 = 40.58
    41 = 7.38
    42 = 34.51
    43 = 195.83
    44 = 286.25
    45 = 286.25
    46 = 95.15000000000001
    47 = 55.18
    48 = 17.2
    49 = 121.74
    50 = 22.14
}

foreach ($r in $pricing.Keys) {
    $ws.Cells.Item($r, 8).Value = $pricing[$r]
}

# TOTAL row (now row 51) - grand total of all line items above
$ws.Range("H51").Value = 4633.380000000001
